# Apply "personalization also with input field" changes to the "To Do" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("To Do")

# --- Update statuses from Open -> Closed (column C) ---
$ws.Range("C6").Value  = "Closed"   # Make weight editable in edit mode
$ws.Range("C13").Value = "Closed"   # When no search selection criteria is selected...

# C15 no longer needs the red "open" highlight font now that it is closed
$ws.Range("C15").ClearFormats()
$ws.Range("C15").Value = "Closed"   # Refactor table as a component and give it the EditMode

$ws.Range("C19").Value = "Closed"   # Unsubscribe change events on Destroy
$ws.Range("C22").Value = "Closed"   # Add charts (fix "Close" typo -> "Closed")

# --- Highlight the still-open personalization related items in red ---
$ws.Range("B23:C23").Font.Color = 255

$ws.Range("C27").Value = "Closed"   # Integrate with NodeServer and Express
$ws.Range("C28").Value = "Closed"   # Study MongoDB schema design (new status)

# --- Fix typo in task description (column B) ---
$ws.Range("B30").Value = "Add a service on the server to call when an error on the client is cought"

# --- New personalization task, inserted right after row 30 ---
$ws.Range("B31").Value = "Centralized error management on the client"
$ws.Range("C31").Value = "Open"

# --- Fix typo in task description (column B) ---
$ws.Range("B42").Value = "how to force reload of a page via router link"
$ws.Range("B42:C42").Font.Color = 255

$ws.Range("C43").Value = "Closed"   # how to format a monetary amount input
$ws.Range("C45").Value = "Closed"   # Which component html5 as free text...

# --- Update the view so the new input field row is the active selection ---
$ws.Activate() | Out-Null
$ws.Range("C6").Select() | Out-Null
